$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the Price (D) column cells we touch so that
# numeric-looking strings (e.g. "330.48", "12.00", "0.00001105") are kept
# as text, matching the original inlineStr/shared-string representation
# instead of being auto-converted to actual numbers by Excel.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.127.03'
$ws.Range('E2').Value = '  -1.46%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.998.35'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  +0.61%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '330.48'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4978'
$ws.Range('E7').Value = '  -0.88%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4194'
$ws.Range('E8').Value = '  -1.65%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '54.81'
$ws.Range('E9').Value = '  +1.87%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.08839'
$ws.Range('E10').Value = '  -3.92%  '
$ws.Range('E11').Value = '  -2.78%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '23.01'
$ws.Range('E12').Value = '  -2.23%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.022.71'
$ws.Range('E13').Value = '  +6.26%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.996'
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.432'
$ws.Range('E15').Value = '  -1.79%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.015'
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '92.51'
$ws.Range('E17').Value = '  -3.52%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001105'
$ws.Range('E18').Value = '  -1.76%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06761'
$ws.Range('E19').Value = '  +1.57%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '19.56'
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.985'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '29.166.67'
$ws.Range('E23').Value = '  -1.32%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '12.00'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.294'
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.246.45'
$ws.Range('E26').Value = '  +4.09%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.85'
$ws.Range('E27').Value = '  +0.36%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '157.21'
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.317'
$ws.Range('E29').Value = '  -3.66%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.258'
$ws.Range('E30').Value = '  -3.51%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '127.20'
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.046'
$ws.Range('E32').Value = '  -0.92%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09866'
$ws.Range('E33').Value = '  -1.03%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.533'
$ws.Range('E34').Value = '  -3.76%  '
$ws.Range('E35').Value = '  -0.66%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.739'
$ws.Range('E36').Value = '  -1.24%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02423'
$ws.Range('E37').Value = '  -2.05%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '9.169'
$ws.Range('E38').Value = '  -4.92%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.314'
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.06394'
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6490'
$ws.Range('E41').Value = '  -1.28%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.59'
$ws.Range('E42').Value = '  -1.84%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1985'
$ws.Range('E43').Value = '  -4.43%  '
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.364'
$ws.Range('E45').Value = '  +4.87%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6212'
$ws.Range('E46').Value = '  -2.28%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '13.37'
$ws.Range('E47').Value = '  -1.54%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.179'
$ws.Range('E48').Value = '  -1.70%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.489'
$ws.Range('E49').Value = '  -1.03%  '
$ws.Range('E50').Value = '  +4.48%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.192'
$ws.Range('E51').Value = '  +9.99%  '
